# Update countries & provincias Spain
# Applies the daily data refresh to the "Pais" worksheet:
#  - updates the "Datos actualizados..." timestamp
#  - updates case counts for several countries
#  - re-sorts a few country pairs whose rank order changed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Julio de 2020 a las 03:21"

# Estados Unidos (row 4) - new totals
$ws.Range("B4").Value = 4371839
$ws.Range("C4").Value = 56130
$ws.Range("D4").Value = 2090129
$ws.Range("E4").Value = 2131861
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 451
$ws.Range("H4").Value = 149849

# India (row 6) - new totals
$ws.Range("B6").Value = 1436019
$ws.Range("C6").Value = 50525
$ws.Range("D6").Value = 918735
$ws.Range("E6").Value = 484472
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 716
$ws.Range("H6").Value = 32812

# Ecuador overtakes Filipinas (rows 31-32)
$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 80694
$ws.Range("C31").Value = 658
$ws.Range("D31").Value = 34896
$ws.Range("E31").Value = 40283
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = 5515

$ws.Range("A32").Value = "Filipinas"
$ws.Range("B32").Value = 80448
$ws.Range("C32").Value = 2036
$ws.Range("D32").Value = 26110
$ws.Range("E32").Value = 52406
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 35
$ws.Range("H32").Value = 1932

# Panama overtakes Emiratos Arabes Unidos (rows 42-43)
$ws.Range("A42").Value = "Panama"
$ws.Range("B42").Value = 60296
$ws.Range("C42").Value = 1432
$ws.Range("D42").Value = 34131
$ws.Range("E42").Value = 24871
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 19
$ws.Range("H42").Value = 1294

$ws.Range("A43").Value = "Emiratos Arabes Unidos"
$ws.Range("B43").Value = 58913
$ws.Range("C43").Value = 351
$ws.Range("D43").Value = 52182
$ws.Range("E43").Value = 6387
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 344

# Kenia overtakes Camerun (rows 68-69)
$ws.Range("A68").Value = "Kenia"
$ws.Range("B68").Value = 17603
$ws.Range("C68").Value = 960
$ws.Range("D68").Value = 7743
$ws.Range("E68").Value = 9580
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 280

$ws.Range("A69").Value = "Camerun"
$ws.Range("B69").Value = 16708
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 14539
$ws.Range("E69").Value = 1784
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 385

# Togo overtakes Santo Tome y Principe (rows 151-152)
$ws.Range("A151").Value = "Togo"
$ws.Range("B151").Value = 868
$ws.Range("C151").Value = 15
$ws.Range("D151").Value = 599
$ws.Range("E151").Value = 251
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 18

$ws.Range("A152").Value = "Santo Tome y Principe"
$ws.Range("B152").Value = 863
$ws.Range("C152").Value = 1
$ws.Range("D152").Value = 696
$ws.Range("E152").Value = 153
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 14

# Guyana overtakes Burundi (rows 164-165)
$ws.Range("A164").Value = "Guyana"
$ws.Range("B164").Value = 370
$ws.Range("C164").Value = 10
$ws.Range("D164").Value = 181
$ws.Range("E164").Value = 169
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 20

$ws.Range("A165").Value = "Burundi"
$ws.Range("B165").Value = 361
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 279
$ws.Range("E165").Value = 81
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 1
